$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to edit cell values, then restore protection afterward.
$ws.Unprotect()

$ws.Range("D2").Value = 0.03134377401633179
$ws.Range("E2").Value = -0.0001732801940738105
$ws.Range("D3").Value = 0.0314970130035889
$ws.Range("E3").Value = 0.02898791708236348
$ws.Range("D4").Value = 0.03154395456677399
$ws.Range("E4").Value = 0.002951666461689761
$ws.Range("D5").Value = 0.06491959996477444
$ws.Range("E5").Value = 0.01106718895193359
$ws.Range("D6").Value = 0.03019933094694324
$ws.Range("E6").Value = -0.001695699090488612
$ws.Range("D7").Value = 0.01566761959333089
$ws.Range("E7").Value = 0.01998217203981567
$ws.Range("D8").Value = 0.03296364589302057
$ws.Range("E8").Value = -0.006508217654570325
$ws.Range("D9").Value = 0.03146190381790088
$ws.Range("E9").Value = 0.002096216329525058
$ws.Range("D10").Value = 0.04740418974788501
$ws.Range("E10").Value = 0.02275098717188029
$ws.Range("D11").Value = 0.02873677149894498
$ws.Range("E11").Value = 0.0108675108675107
$ws.Range("D12").Value = 0.01540711331499379
$ws.Range("E12").Value = 0.02758438353749892
$ws.Range("D13").Value = 0.01713464042956574
$ws.Range("E13").Value = -0.01171676006113098
$ws.Range("D14").Value = 0.01440640453618438
$ws.Range("E14").Value = 0.009519321394910518
$ws.Range("D15").Value = 0.007166347321789848
$ws.Range("E15").Value = -0.01997563946406833
$ws.Range("D16").Value = 0.007661173469249202
$ws.Range("E16").Value = -0.01478630747417453
$ws.Range("D17").Value = 0.03262962369548039
$ws.Range("E17").Value = 0.0247596854063501
$ws.Range("D18").Value = 0.02919338488629085
$ws.Range("E18").Value = 0.01615260926765094
$ws.Range("D19").Value = 0.032166609185882
$ws.Range("E19").Value = -0.003799071338117299
$ws.Range("D20").Value = 0.03212781450556374
$ws.Range("E20").Value = 0.01663346012195865
$ws.Range("D21").Value = 0.0487635553462367
$ws.Range("E21").Value = 0.005214942341274398
$ws.Range("D22").Value = 0.02933265778863339
$ws.Range("E22").Value = -0.01356963364634312
$ws.Range("D23").Value = 0.03014055700626107
$ws.Range("E23").Value = -0.004665830035074148
$ws.Range("D24").Value = 0.02884869415166315
$ws.Range("E24").Value = -0.003872919818456722
$ws.Range("D25").Value = 0.01398005099948675
$ws.Range("E25").Value = -0.0237540237540238
$ws.Range("D26").Value = 0.01480055848821786
$ws.Range("E26").Value = -0.0174569473932531
$ws.Range("D27").Value = 0.03227872581200175
$ws.Range("E27").Value = -0.008262823902696947
$ws.Range("D28").Value = 0.03033744000887622
$ws.Range("E28").Value = -0.01324808184143222
$ws.Range("D29").Value = 0.03001428032182515
$ws.Range("E29").Value = 0.007755244484082269
$ws.Range("D30").Value = 0.02815213566654886
$ws.Range("E30").Value = -0.004395937547369888
$ws.Range("D31").Value = 0.02766099501371974
$ws.Range("E31").Value = 0.01085538772247241
$ws.Range("D32").Value = 0.02864793168101618
$ws.Range("E32").Value = 0.005958426433746222
$ws.Range("D33").Value = 0.02981856115961955
$ws.Range("E33").Value = 0.02317124735729381
$ws.Range("D34").Value = 0.03087203070366179
$ws.Range("E34").Value = 0.003191836939857806
$ws.Range("D35").Value = 0.0310919965410663
$ws.Range("E35").Value = -0.01696924324661553
$ws.Range("D36").Value = 0.03162891491667096
$ws.Range("E36").Value = 0.0005580836266849687
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 0.004224158766453101

# Update the confidential disclosure footnote date (2021-06-10 -> 2021-06-14).
# Use Cells.Replace (find/replace) rather than re-assigning .Value so the
# shared-string table is updated cleanly instead of leaving a stale,
# unreferenced string entry behind.
$ws.Cells.Replace("2021-06-10", "2021-06-14")

# Restore sheet protection.
$ws.Protect()
